# Harvard case classification update: "_old" app scoring columns were
# recalculated, and the average_doctor / average_doctor_old header order
# was swapped (average_doctor_old now precedes average_doctor).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header swap
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4
$ws.Range("E4").Value = 0.425
$ws.Range("F4").Value = 0.07199999999999999
$ws.Range("G4").Value = 0.268
$ws.Range("N4").Value = 0.433
$ws.Range("O4").Value = 0.065
$ws.Range("P4").Value = 0.255
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.129
$ws.Range("W4").Value = 0.282
$ws.Range("X4").Value = 0.11
$ws.Range("Y4").Value = 0.331
$ws.Range("AI4").Value = 0.284
$ws.Range("AJ4").Value = 0.08699999999999999
$ws.Range("AK4").Value = 0.296
$ws.Range("AU4").Value = 0.19
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.167
$ws.Range("BA4").Value = 2.001
$ws.Range("BB4").Value = 0.158
$ws.Range("BC4").Value = 0.397
$ws.Range("BG4").Value = 0.729
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.376
$ws.Range("BM4").Value = 0.716
$ws.Range("BN4").Value = 0.08
$ws.Range("BO4").Value = 0.282
$ws.Range("BP4").Value = 0.667
$ws.Range("BQ4").Value = 0.708

# Row 5
$ws.Range("E5").Value = 0.539
$ws.Range("F5").Value = 0.08599999999999999
$ws.Range("G5").Value = 0.294
$ws.Range("N5").Value = 0.735
$ws.Range("O5").Value = 0.081
$ws.Range("P5").Value = 0.285
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.083
$ws.Range("W5").Value = 0.27
$ws.Range("X5").Value = 0.109
$ws.Range("Y5").Value = 0.331
$ws.Range("AI5").Value = 0.302
$ws.Range("AJ5").Value = 0.094
$ws.Range("AK5").Value = 0.307
$ws.Range("AU5").Value = 0.369
$ws.Range("AV5").Value = 0.093
$ws.Range("AW5").Value = 0.305
$ws.Range("BA5").Value = 1.335
$ws.Range("BB5").Value = 0.079
$ws.Range("BC5").Value = 0.282
$ws.Range("BG5").Value = 0.397
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.226
$ws.Range("BM5").Value = 0.549
$ws.Range("BN5").Value = 0.064
$ws.Range("BO5").Value = 0.253
$ws.Range("BP5").Value = 0.445
$ws.Range("BQ5").Value = 0.455

# Row 6
$ws.Range("E6").Value = 0.475
$ws.Range("N6").Value = 0.545
$ws.Range("Q6").Value = 0.019
$ws.Range("W6").Value = 0.276
$ws.Range("AI6").Value = 0.293
$ws.Range("AU6").Value = 0.251
$ws.Range("BA6").Value = 1.593
$ws.Range("BG6").Value = 0.514
$ws.Range("BM6").Value = 0.621
$ws.Range("BP6").Value = 0.531
$ws.Range("BQ6").Value = 0.551

# Row 7
$ws.Range("E7").Value = 0.512
$ws.Range("N7").Value = 0.645
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.272
$ws.Range("AI7").Value = 0.298
$ws.Range("AU7").Value = 0.31
$ws.Range("BA7").Value = 1.427
$ws.Range("BG7").Value = 0.437
$ws.Range("BM7").Value = 0.576
$ws.Range("BP7").Value = 0.476
$ws.Range("BQ7").Value = 0.489

# Row 8
$ws.Range("E8").Value = 0.602
$ws.Range("F8").Value = 0.112
$ws.Range("G8").Value = 0.335
$ws.Range("N8").Value = 0.772
$ws.Range("O8").Value = 0.066
$ws.Range("P8").Value = 0.258
$ws.Range("Q8").Value = 0.017
$ws.Range("W8").Value = 0.295
$ws.Range("X8").Value = 0.118
$ws.Range("Y8").Value = 0.344
$ws.Range("AI8").Value = 0.324
$ws.Range("AJ8").Value = 0.126
$ws.Range("AK8").Value = 0.355
$ws.Range("AU8").Value = 0.311
$ws.Range("AV8").Value = 0.08400000000000001
$ws.Range("AW8").Value = 0.291
$ws.Range("BA8").Value = 1.741
$ws.Range("BB8").Value = 0.124
$ws.Range("BC8").Value = 0.352
$ws.Range("BG8").Value = 0.5649999999999999
$ws.Range("BH8").Value = 0.107
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.6889999999999999
$ws.Range("BN8").Value = 0.067
$ws.Range("BO8").Value = 0.259
$ws.Range("BP8").Value = 0.58
$ws.Range("BQ8").Value = 0.603

# Row 9
$ws.Range("E9").Value = 0.548
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.677
$ws.Range("O9").Value = 0.219
$ws.Range("P9").Value = 0.467
$ws.Range("W9").Value = 0.194
$ws.Range("X9").Value = 0.156
$ws.Range("Y9").Value = 0.395
$ws.Range("AI9").Value = 0.247
$ws.Range("AJ9").Value = 0.186
$ws.Range("AK9").Value = 0.431
$ws.Range("BA9").Value = 1.699
$ws.Range("BB9").Value = 0.248
$ws.Range("BC9").Value = 0.498
$ws.Range("BG9").Value = 0.602
$ws.Range("BH9").Value = 0.24
$ws.Range("BI9").Value = 0.489
$ws.Range("BM9").Value = 0.645
$ws.Range("BN9").Value = 0.229
$ws.Range("BO9").Value = 0.478
$ws.Range("BP9").Value = 0.5659999999999999
$ws.Range("BQ9").Value = 0.585

# Row 10
$ws.Range("E10").Value = 0.677
$ws.Range("F10").Value = 0.219
$ws.Range("G10").Value = 0.467
$ws.Range("N10").Value = 0.871
$ws.Range("O10").Value = 0.112
$ws.Range("P10").Value = 0.335
$ws.Range("W10").Value = 0.366
$ws.Range("X10").Value = 0.232
$ws.Range("Y10").Value = 0.482
$ws.Range("AI10").Value = 0.355
$ws.Range("AJ10").Value = 0.229
$ws.Range("AK10").Value = 0.478
$ws.Range("AU10").Value = 0.301
$ws.Range("AV10").Value = 0.21
$ws.Range("AW10").Value = 0.459
$ws.Range("BA10").Value = 2.076
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.656
$ws.Range("BH10").Value = 0.226
$ws.Range("BI10").Value = 0.475
$ws.Range("BM10").Value = 0.839
$ws.Range("BN10").Value = 0.135
$ws.Range("BO10").Value = 0.368
$ws.Range("BP10").Value = 0.6919999999999999
$ws.Range("BQ10").Value = 0.723

# Row 11
$ws.Range("E11").Value = 0.71
$ws.Range("F11").Value = 0.206
$ws.Range("G11").Value = 0.454
$ws.Range("N11").Value = 0.892
$ws.Range("O11").Value = 0.096
$ws.Range("P11").Value = 0.31
$ws.Range("W11").Value = 0.366
$ws.Range("X11").Value = 0.232
$ws.Range("Y11").Value = 0.482
$ws.Range("AI11").Value = 0.387
$ws.Range("AJ11").Value = 0.237
$ws.Range("AK11").Value = 0.487
$ws.Range("AU11").Value = 0.441
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.076
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.656
$ws.Range("BH11").Value = 0.226
$ws.Range("BI11").Value = 0.475
$ws.Range("BM11").Value = 0.839
$ws.Range("BN11").Value = 0.135
$ws.Range("BO11").Value = 0.368
$ws.Range("BP11").Value = 0.6919999999999999
$ws.Range("BQ11").Value = 0.726

# Row 12
$ws.Range("E12").Value = 1.409
$ws.Range("F12").Value = 0.757
$ws.Range("G12").Value = 0.87
$ws.Range("N12").Value = 1.471
$ws.Range("O12").Value = 1.049
$ws.Range("P12").Value = 1.024
$ws.Range("W12").Value = 1.647
$ws.Range("X12").Value = 0.581
$ws.Range("Y12").Value = 0.762
$ws.Range("AI12").Value = 1.722
$ws.Range("AJ12").Value = 1.312
$ws.Range("AK12").Value = 1.145
$ws.Range("AU12").Value = 2.767
$ws.Range("AV12").Value = 2.737
$ws.Range("AW12").Value = 1.654
$ws.Range("BA12").Value = 3.708
$ws.Range("BB12").Value = 0.401
$ws.Range("BC12").Value = 0.633
$ws.Range("BG12").Value = 1.098
$ws.Range("BH12").Value = 0.121
$ws.Range("BI12").Value = 0.349
$ws.Range("BM12").Value = 1.295
$ws.Range("BN12").Value = 0.336
$ws.Range("BO12").Value = 0.58
$ws.Range("BP12").Value = 1.236
$ws.Range("BQ12").Value = 1.261

# Row 13
$ws.Range("E13").Value = 1.563
$ws.Range("F13").Value = 0.65
$ws.Range("G13").Value = 0.806
$ws.Range("N13").Value = 2.074
$ws.Range("O13").Value = 0.9350000000000001
$ws.Range("P13").Value = 0.967
$ws.Range("W13").Value = 1.031
$ws.Range("X13").Value = 0.193
$ws.Range("Y13").Value = 0.439
$ws.Range("AI13").Value = 1.277
$ws.Range("AJ13").Value = 0.374
$ws.Range("AK13").Value = 0.612
$ws.Range("AU13").Value = 2.285
$ws.Range("AV13").Value = 0.925
$ws.Range("AW13").Value = 0.962
$ws.Range("BA13").Value = 2.353
$ws.Range("BB13").Value = 0.297
$ws.Range("BC13").Value = 0.545
$ws.Range("BG13").Value = 0.585
$ws.Range("BH13").Value = 0.07099999999999999
$ws.Range("BI13").Value = 0.267
$ws.Range("BM13").Value = 0.896
$ws.Range("BN13").Value = 0.283
$ws.Range("BO13").Value = 0.532
$ws.Range("BP13").Value = 0.784
$ws.Range("BQ13").Value = 0.727
